# Fruta / hortaliza, semanal
# A new weekly record was added to the data set. This is modeled as inserting
# a new row at position 904 (pushing the existing rows 904-992 down to
# 905-993) and then populating the new row 904 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 904; this shifts rows 904:992 down
# to 905:993 (and automatically extends the used range / dimension).
$ws.Rows("904:904").Insert()

# Populate the newly inserted row 904 with the new record.
$ws.Cells.Item(904, 1).Value = 4
$ws.Cells.Item(904, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(904, 3).Value = "Los Lagos"
$ws.Cells.Item(904, 4).Value = 45166
$ws.Cells.Item(904, 5).Value = 10
$ws.Cells.Item(904, 6).Value = 100112004
$ws.Cells.Item(904, 7).Value = "Cebolla"
$ws.Cells.Item(904, 8).Value = "Sin especificar"
$ws.Cells.Item(904, 9).Value = "1a (guarda)"
$ws.Cells.Item(904, 10).Value = 300
$ws.Cells.Item(904, 11).Value = 16000
$ws.Cells.Item(904, 12).Value = 16000
$ws.Cells.Item(904, 13).Value = 16000
$ws.Cells.Item(904, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(904, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(904, 16).Value = 889
$ws.Cells.Item(904, 17).Value = 18
$ws.Cells.Item(904, 18).Value = "Hortaliza"
